$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Code" column (A) text duplicates with a sequential numeric id
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(4, 1).Value = 3

# Widen column A to fit the new content
$ws.Columns.Item(1).ColumnWidth = 18.140625

# Reflect the new active cell selection
$ws.Range("D4").Select()
